$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the missing role/batch/semester/branch ids for existing rows 4 & 5 ---
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 6
$ws.Range("G4").Value = 1

$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 6
$ws.Range("G5").Value = 5

# --- New row 6: Ameya Bhattacharya ---
$ws.Range("A6").Value = 5678
$ws.Range("B6").Value = "Ameya Bhattacharya"
$ws.Range("C6").Value = "ameya.bhattacharya@sitpune.edu.in"
$ws.Range("D6").Value = 3

# --- New row 7: Aishwarya Singh ---
$ws.Range("A7").Value = 102
$ws.Range("B7").Value = "Aishwarya Singh"
$ws.Range("C7").Value = "1994aishwaryasingh@gmail.com"
$ws.Range("D7").Value = 2

# --- Email hyperlinks for the two new rows ---
$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:ameya.bhattacharya@sitpune.edu.in", "", "", "ameya.bhattacharya@sitpune.edu.in")
$ws.Hyperlinks.Add($ws.Range("C7"), "mailto:1994aishwaryasingh@gmail.com", "", "", "1994aishwaryasingh@gmail.com")

# Hyperlinks.Add stamps its own default "Hyperlink" style (themed blue + underline);
# restore the same look already used by the existing email cells (C2:C5).
$ws.Range("C6").Font.Name = $ws.Range("C4").Font.Name
$ws.Range("C6").Font.Size = $ws.Range("C4").Font.Size
$ws.Range("C6").Font.Color = $ws.Range("C4").Font.Color
$ws.Range("C6").Font.Underline = $false

$ws.Range("C7").Font.Name = $ws.Range("C4").Font.Name
$ws.Range("C7").Font.Size = $ws.Range("C4").Font.Size
$ws.Range("C7").Font.Color = $ws.Range("C4").Font.Color
$ws.Range("C7").Font.Underline = $false

# --- Widen column C slightly to fit the new, longer email addresses ---
$ws.Columns.Item(3).ColumnWidth = 31.333333333333332

# --- Move the active selection, matching the author's last editing position ---
$ws.Range("E8").Select()
